$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 371, pushing existing rows 371..464 down to 372..465.
$ws.Rows.Item(371).Insert()

# Populate the new row 371 with its data. Columns A,B,C,E,F,G,H,I,J,L,Q,R,T
# keep the same values the (now shifted) row carries in this block, only
# D, K, M, N, O, P, S differ per the target data set.
$ws.Range("A371").Value = 5
$ws.Range("B371").Value = "Macroferia Regional de Talca"
$ws.Range("C371").Value = "Maule"
$ws.Range("D371").Value = 45204
$ws.Range("E371").Value = 7
$ws.Range("F371").Value = "Fruta"
$ws.Range("G371").Value = 100108
$ws.Range("H371").Value = "Tropicales y subtropicales"
$ws.Range("I371").Value = 100108005
$ws.Range("J371").Value = "Piña"
$ws.Range("K371").Value = "Sin especificar"
$ws.Range("L371").Value = "Segunda"
$ws.Range("M371").Value = 200
$ws.Range("N371").Value = 21000
$ws.Range("O371").Value = 21000
$ws.Range("P371").Value = 21000
$ws.Range("Q371").Value = "$/caja 14 unidades"
$ws.Range("R371").Value = "Ecuador"
$ws.Range("S371").Value = 1500
$ws.Range("T371").Value = 14

# Excel's Date format style (s="2") for the date column, matching the other
# rows in this block.
$ws.Range("D371").NumberFormat = $ws.Range("D372").NumberFormat
